$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 3.4
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 1.73
$ws.Range("T2").Value = 2.08
